$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Figure out the next empty row right after the current localization data
$lastRow = $ws.UsedRange.Rows.Count
$newRow = $lastRow + 1

# Copy the formatting from a row that uses the standard data style in both
# columns (the immediately preceding row has a one-off style in column B)
# onto the new row before filling in the values, so the new cells end up
# sharing the same cell style as the rest of the table.
$formatRow = $lastRow - 1
$ws.Range("A" + $formatRow + ":B" + $formatRow).Copy()
$ws.Range("A" + $newRow + ":B" + $newRow).PasteSpecial(-4122)  # xlPasteFormats

# Add the missing localization key/value pair
$ws.Cells.Item($newRow, 1).Value = "Warning"
$ws.Cells.Item($newRow, 2).Value = "Uyari"
